# Update "想去人数" (interest count) figures in column F across the
# workbook's sheets, reflecting refreshed scrape data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 377
$ws1.Cells.Item(5, 6).Value = 1286
$ws1.Cells.Item(7, 6).Value = 2464
$ws1.Cells.Item(8, 6).Value = 871
$ws1.Cells.Item(9, 6).Value = 18510
$ws1.Cells.Item(10, 6).Value = 49
$ws1.Cells.Item(11, 6).Value = 1873
$ws1.Cells.Item(12, 6).Value = 654
$ws1.Cells.Item(14, 6).Value = 319
$ws1.Cells.Item(16, 6).Value = 193
$ws1.Cells.Item(18, 6).Value = 66
$ws1.Cells.Item(19, 6).Value = 316
$ws1.Cells.Item(23, 6).Value = 76

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value = 28
$ws2.Cells.Item(8, 6).Value = 120
$ws2.Cells.Item(9, 6).Value = 110
$ws2.Cells.Item(14, 6).Value = 67

# Sheet "本地生活" (index 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 5862
$ws3.Cells.Item(3, 6).Value = 554
$ws3.Cells.Item(4, 6).Value = 551

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 5862
$ws4.Cells.Item(4, 6).Value = 554
$ws4.Cells.Item(5, 6).Value = 551
$ws4.Cells.Item(6, 6).Value = 377
$ws4.Cells.Item(10, 6).Value = 1286
$ws4.Cells.Item(14, 6).Value = 28
$ws4.Cells.Item(15, 6).Value = 2464
$ws4.Cells.Item(16, 6).Value = 871
$ws4.Cells.Item(17, 6).Value = 18510
$ws4.Cells.Item(18, 6).Value = 49
$ws4.Cells.Item(20, 6).Value = 120
$ws4.Cells.Item(21, 6).Value = 120
$ws4.Cells.Item(22, 6).Value = 1873
$ws4.Cells.Item(23, 6).Value = 654
$ws4.Cells.Item(24, 6).Value = 110
$ws4.Cells.Item(26, 6).Value = 319
$ws4.Cells.Item(28, 6).Value = 193
$ws4.Cells.Item(31, 6).Value = 66
$ws4.Cells.Item(34, 6).Value = 316
$ws4.Cells.Item(36, 6).Value = 67
$ws4.Cells.Item(48, 6).Value = 76

Write-Host "F column updates applied"
